$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '39.847.99'
$ws.Range("E2").Value = '  +0.27%  '

# Row 3
$ws.Range("D3").Value = '2.218.77'
$ws.Range("E3").Value = '  +0.23%  '

# Row 4
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '291.75'
$ws.Range("E5").Value = '  -0.20%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '86.85'
$ws.Range("E6").Value = '  +0.28%  '

# Row 7
$ws.Range("E7").Value = '  -0.58%  '

# Row 8
$ws.Range("E8").Value = '  +0.10%  '

# Row 9
$ws.Range("E9").Value = '  -0.84%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.39'
$ws.Range("E10").Value = '  +0.57%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0779'
$ws.Range("E11").Value = '  -0.74%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '49.96'
$ws.Range("E12").Value = '  +5.50%  '

# Row 13
$ws.Range("E13").Value = '  +2.65%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.43'
$ws.Range("E14").Value = '  +1.39%  '

# Row 15
$ws.Range("D15").Value = '2.565.72'
$ws.Range("E15").Value = '  +0.54%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.78'
$ws.Range("E16").Value = '  -1.81%  '

# Row 17
$ws.Range("D17").Value = '2.227.40'
$ws.Range("E17").Value = '  +1.02%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.730'
$ws.Range("E18").Value = '  +0.35%  '

# Row 19
$ws.Range("D19").Value = '39.812.89'
$ws.Range("E19").Value = '  +0.32%  '

# Row 20
$ws.Range("E20").Value = '  +0.61%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.05'
$ws.Range("E21").Value = '  -4.12%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.74'
$ws.Range("E22").Value = '  -1.07%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.80'
$ws.Range("E23").Value = '  -0.09%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '237.61'
$ws.Range("E24").Value = '  +0.89%  '

# Row 25
$ws.Range("E25").Value = '  +0.10%  '

# Row 26
$ws.Range("E26").Value = '  -0.36%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.82'
$ws.Range("E27").Value = '  -0.24%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.08'
$ws.Range("E28").Value = '  +1.56%  '

# Row 29
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.23'
$ws.Range("E29").Value = '  -0.35%  '

# Row 30
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.03'
$ws.Range("E30").Value = '  -7.60%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '155.94'
$ws.Range("E31").Value = '  +2.60%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.88'
$ws.Range("E32").Value = '  -2.68%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.18%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.96'
$ws.Range("E34").Value = '  +0.61%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.98'
$ws.Range("E35").Value = '  +6.68%  '

# Row 36
$ws.Range("E36").Value = '  -0.37%  '

# Row 37
$ws.Range("E37").Value = '  -1.89%  '

# Row 38
$ws.Range("E38").Value = '  -0.38%  '

# Row 39
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0993'
$ws.Range("E39").Value = '  +0.76%  '

# Row 40
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.74'
$ws.Range("E40").Value = '  +2.78%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.30'
$ws.Range("E41").Value = '  -4.48%  '

# Row 42
$ws.Range("D42").Value = '2.104.49'
$ws.Range("E42").Value = '  +1.18%  '

# Row 43
$ws.Range("E43").Value = '  -1.98%  '

# Row 44
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0271'
$ws.Range("E44").Value = '  +0.99%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.03'
$ws.Range("E45").Value = '  +2.14%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.87'
$ws.Range("E46").Value = '  -1.33%  '

# Row 47
$ws.Range("E47").Value = '  -7.76%  '

# Row 48
$ws.Range("E48").Value = '  +4.03%  '

# Row 49
$ws.Range("D49").Value = '2.439.81'
$ws.Range("E49").Value = '  +0.58%  '

# Row 50
$ws.Range("E50").Value = '  +1.18%  '

# Row 51
$ws.Range("E51").Value = '  +2.40%  '
